$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# PlantMaster table is changing from a "space_needed / sun" model to a
# "depth_needed / capacity_needed" model (container-allocation columns),
# and two new rows (containers / plants) are being appended: Passion
# Fruit and Cotton.
# ----------------------------------------------------------------------

# --- Header row --------------------------------------------------------
# A1 ("name") stays the same.
$ws.Range("B1").Value = "depth_needed"

# Row 10 (Cotton) and row 9 (Passion Fruit) are written before the C1
# header so the workbook's shared-string table ends up in the same order
# as the source file (depth_needed, Cotton, capacity_needed,
# Passion Fruit).
$ws.Cells.Item(10, 1).Value = "Cotton"
$ws.Cells.Item(10, 2).Value = 90
$ws.Cells.Item(10, 3).Value = 15

$ws.Range("C1").Value = "capacity_needed"

$ws.Cells.Item(9, 1).Value = "Passion Fruit"
$ws.Cells.Item(9, 2).Value = 30
$ws.Cells.Item(9, 3).Value = 20

# --- Updated depth_needed / capacity_needed values for existing plants -
$ws.Cells.Item(2, 2).Value = 20   # Basil
$ws.Cells.Item(2, 3).Value = 4

$ws.Cells.Item(3, 2).Value = 15   # Bush Bean
$ws.Cells.Item(3, 3).Value = 8

$ws.Cells.Item(4, 2).Value = 15   # Chive
$ws.Cells.Item(4, 3).Value = 2

$ws.Cells.Item(5, 2).Value = 20   # Garlic
$ws.Cells.Item(5, 3).Value = 5

$ws.Cells.Item(6, 2).Value = 14   # Radish
$ws.Cells.Item(6, 3).Value = 1

$ws.Cells.Item(7, 2).Value = 30   # Stevia
$ws.Cells.Item(7, 3).Value = 10

$ws.Cells.Item(8, 2).Value = 30   # Zucchini
$ws.Cells.Item(8, 3).Value = 15

# --- Column widths (re-fitted to the new, wider header/content) --------
$ws.Columns.Item(1).ColumnWidth = 11.5
$ws.Columns.Item(2).ColumnWidth = 12.333333333333332
$ws.Columns.Item(3).ColumnWidth = 14.5

# --- Selection moves to A10, matching the new bottom of the table ------
$ws.Range("A10").Select()
